$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 696, pushing the existing rows
# (old 696..763) down to 698..765.
$ws.Rows.Item(696).EntireRow.Insert()
$ws.Rows.Item(696).EntireRow.Insert()

# New record #1 (row 696) - Tomate, Larga vida, Primera
$ws.Cells.Item(696, 1).Value = 7
$ws.Cells.Item(696, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(696, 3).Value = "Ñuble"
$ws.Cells.Item(696, 4).Value = 45106
$ws.Cells.Item(696, 5).Value = 16
$ws.Cells.Item(696, 6).Value = 100112020
$ws.Cells.Item(696, 7).Value = "Tomate"
$ws.Cells.Item(696, 8).Value = "Larga vida"
$ws.Cells.Item(696, 9).Value = "Primera"
$ws.Cells.Item(696, 10).Value = 100
$ws.Cells.Item(696, 11).Value = 13000
$ws.Cells.Item(696, 12).Value = 13000
$ws.Cells.Item(696, 13).Value = 13000
$ws.Cells.Item(696, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(696, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(696, 16).Value = 722
$ws.Cells.Item(696, 17).Value = 18
$ws.Cells.Item(696, 18).Value = "Hortaliza"

# New record #2 (row 697) - Tomate, Larga vida, Segunda
$ws.Cells.Item(697, 1).Value = 7
$ws.Cells.Item(697, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(697, 3).Value = "Ñuble"
$ws.Cells.Item(697, 4).Value = 45106
$ws.Cells.Item(697, 5).Value = 16
$ws.Cells.Item(697, 6).Value = 100112020
$ws.Cells.Item(697, 7).Value = "Tomate"
$ws.Cells.Item(697, 8).Value = "Larga vida"
$ws.Cells.Item(697, 9).Value = "Segunda"
$ws.Cells.Item(697, 10).Value = 120
$ws.Cells.Item(697, 11).Value = 11000
$ws.Cells.Item(697, 12).Value = 11000
$ws.Cells.Item(697, 13).Value = 11000
$ws.Cells.Item(697, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(697, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(697, 16).Value = 611
$ws.Cells.Item(697, 17).Value = 18
$ws.Cells.Item(697, 18).Value = "Hortaliza"
